$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header cell C6: "Time" -> "Time (h)"
$ws.Range("C6").Value = "Time (h)"

# 2. New row 9 date cell first, so the date-format style is registered
#    before the wrap-text style (keeps style index ordering stable).
$ws.Range("B9").Value = 45348
$ws.Range("B9").NumberFormat = "d-mmm"

# 3. Apply wrap-text formatting to D6:E8 (existing data rows)
$ws.Range("D6:E8").WrapText = $true | Out-Null

# 4. Remaining new row 9 data
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = "Started bug fixing and writing more functionality to fit my requirtements. Also fixed all of the ship prefabs and implemented a simple UI system."
$ws.Range("E9").Value = "Placing ships still glitchy. Think its to do with colliders."
$ws.Range("D9:E9").WrapText = $true | Out-Null
$ws.Rows.Item(9).RowHeight = 28

# 5. Select E9 last (to match final selection state)
$ws.Range("E9").Select() | Out-Null

Write-Host "done"
